$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking decimal string.
# These must be forced to Text format first, otherwise Excel would
# coerce them (e.g. "1.00" -> 1) and the original fixed-decimal
# text representation used throughout this price column would be lost.
$numericTextCells = [ordered]@{}
$numericTextCells['D5'] = '563.26'
$numericTextCells['D6'] = '142.59'
$numericTextCells['D11'] = '5.23'
$numericTextCells['D13'] = '25.60'
$numericTextCells['D18'] = '11.25'
$numericTextCells['D19'] = '321.61'
$numericTextCells['D20'] = '4.14'
$numericTextCells['D21'] = '6.83'
$numericTextCells['D23'] = '65.99'
$numericTextCells['D25'] = '8.66'
$numericTextCells['D26'] = '564.03'
$numericTextCells['D30'] = '8.19'
$numericTextCells['D31'] = '1.40'
$numericTextCells['D35'] = '1.00'
$numericTextCells['D36'] = '4.78'
$numericTextCells['D38'] = '152.67'
$numericTextCells['D40'] = '18.54'
$numericTextCells['D41'] = '1.80'
$numericTextCells['D42'] = '1.00'
$numericTextCells['D44'] = '2.24'
$numericTextCells['D47'] = '19.89'
$numericTextCells['D48'] = '0.592'

foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
}

# Remaining cells (text, urls, percentages, coin names, multi-dot
# formatted big numbers) are not valid numeric literals, so plain
# assignment keeps them as text without any extra formatting work.
$ws.Range('D2').Value = '62.025.59'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.417.13'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '2.848.47'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '61.876.18'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '2.412.97'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('E18').Value = '  +1.43%  '
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('E25').Value = '  -4.93%  '
$ws.Range('E26').Value = '  -2.14%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '2.525.78'
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('E31').Value = '  -3.52%  '
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('E38').Value = '  +2.71%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('E41').Value = '  -3.51%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('E51').Value = '  +0.60%  '

